# Actualización SmartScore desde Streamlit (Paula Belén Chairez Rosas)
# Adds a new data row (row 12) to Sheet1 with the participant's results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

# --- A: ID_Participante ---
$ws.Cells.Item($row, 1).Value = "Paula Belén Chairez Rosas_20251120_205520"

# --- B: Grupo_Experimental (empty for this participant) ---
# Force a text-typed, empty-string cell (not a blank/missing cell).
$ws.Cells.Item($row, 2).Value = "'"

# --- C: Nombre Completo ---
$ws.Cells.Item($row, 3).Value = "Paula Belén Chairez Rosas"

# --- D: Edad (numeric) ---
$ws.Cells.Item($row, 4).Value = 20

# --- E: Género ---
$ws.Cells.Item($row, 5).Value = "Femenino"

# --- F: Fecha ---
$ws.Cells.Item($row, 6).Value = "2025-11-20 20:55:20"

# --- G: Pesos (multi-line JSON, stored as text) ---
$pesos = @"
{
  "portion": 0.0,
  "diet": 0.0,
  "salt": 0.0,
  "fat": 0.0,
  "natural": 0.0,
  "convenience": 0.0,
  "price": 0.0
}
"@
$ws.Cells.Item($row, 7).Value = $pesos

# --- Instant Noodles Top 1 (H, I, J) ---
$ws.Cells.Item($row, 8).Value = "Maruchan Ramen Sabor Pollo"
$ws.Cells.Item($row, 9).Value = "'0.000"
$ws.Cells.Item($row, 10).Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"

# --- Instant Noodles Top 2 (K, L, M) ---
$ws.Cells.Item($row, 11).Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Cells.Item($row, 12).Value = "'0.000"
$ws.Cells.Item($row, 13).Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

# --- Instant Noodles Top 3 (N, O, P) ---
$ws.Cells.Item($row, 14).Value = "Nongshim Shin Ramyun"
$ws.Cells.Item($row, 15).Value = "'0.000"
$ws.Cells.Item($row, 16).Value = "Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio"

# --- Mac & Cheese Top 1 (Q, R, S) ---
$ws.Cells.Item($row, 17).Value = "Annie’s Shells & White Cheddar"
$ws.Cells.Item($row, 18).Value = "'0.000"
$ws.Cells.Item($row, 19).Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"

# --- Mac & Cheese Top 2 (T, U, V) ---
$ws.Cells.Item($row, 20).Value = "Velveeta Original Shells & Cheese (microwave cups)"
$ws.Cells.Item($row, 21).Value = "'0.000"
$ws.Cells.Item($row, 22).Value = "Muy cremoso, porción individual, rápido, salado, ideal para niños"

# --- Mac & Cheese Top 3 (W, X, Y) ---
$ws.Cells.Item($row, 23).Value = "Kraft Macaroni & Cheese Dinner"
$ws.Cells.Item($row, 24).Value = "'0.000"
$ws.Cells.Item($row, 25).Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

# --- Top 1 (Z, AA, AB) ---
$ws.Cells.Item($row, 26).Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Cells.Item($row, 27).Value = "'0.000"
$ws.Cells.Item($row, 28).Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

# --- Top 2 (AC, AD, AE) ---
$ws.Cells.Item($row, 29).Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Cells.Item($row, 30).Value = "'0.000"
$ws.Cells.Item($row, 31).Value = "Portátil, saludable, fácil, buena textura, sabor suave"

# --- Top 3 (AF, AG, AH) ---
$ws.Cells.Item($row, 32).Value = "Kitchens of India Variety Pack"
$ws.Cells.Item($row, 33).Value = "'0.000"
$ws.Cells.Item($row, 34).Value = "Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad"

# Re-fit the row height for the multi-line "Pesos" cell BEFORE stripping the
# quote-prefix formatting below (AutoFit can re-touch cell formatting, so it
# must run first or it will re-introduce styles we just cleared).
$ws.Rows.Item($row).AutoFit()

# Strip the quote-prefix formatting introduced by forcing text values above,
# restoring these cells to the default (unstyled) appearance.
$ws.Cells.Item($row, 2).ClearFormats()
$ws.Cells.Item($row, 9).ClearFormats()
$ws.Cells.Item($row, 12).ClearFormats()
$ws.Cells.Item($row, 15).ClearFormats()
$ws.Cells.Item($row, 18).ClearFormats()
$ws.Cells.Item($row, 21).ClearFormats()
$ws.Cells.Item($row, 24).ClearFormats()
$ws.Cells.Item($row, 27).ClearFormats()
$ws.Cells.Item($row, 30).ClearFormats()
$ws.Cells.Item($row, 33).ClearFormats()
